$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that are formatted/scraped as text (they
# frequently contain thousands-separator dots, trailing zeros, or very small
# decimals). Force the "Price" column to a Text number format before writing
# so Excel does not silently reinterpret the values as numbers (which would
# drop trailing zeros / switch to scientific notation).
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.978.00'
$ws.Range("E2").Value = '  -2.98%  '
$ws.Range("D3").Value = '1.859.58'
$ws.Range("E3").Value = '  -2.50%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '306.19'
$ws.Range("E5").Value = '  -2.19%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.5106'
$ws.Range("E7").Value = '  +2.54%  '
$ws.Range("D8").Value = '0.3736'
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("D9").Value = '0.07098'
$ws.Range("E9").Value = '  -2.23%  '
$ws.Range("D10").Value = '0.8880'
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").Value = '20.53'
$ws.Range("E11").Value = '  -2.76%  '
$ws.Range("D12").Value = '0.07541'
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("D13").Value = '1.851.63'
$ws.Range("E13").Value = '  -3.96%  '
$ws.Range("D14").Value = '5.288'
$ws.Range("E14").Value = '  -3.16%  '
$ws.Range("D15").Value = '88.48'
$ws.Range("E15").Value = '  -3.70%  '
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '0.000008387'
$ws.Range("E17").Value = '  -3.59%  '
$ws.Range("D18").Value = '14.05'
$ws.Range("E18").Value = '  -3.34%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").Value = '27.035.49'
$ws.Range("E20").Value = '  -2.89%  '
$ws.Range("D21").Value = '5.054'
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").Value = '2.097.61'
$ws.Range("E22").Value = '  -1.52%  '
$ws.Range("E23").Value = '  -2.72%  '
$ws.Range("D24").Value = '6.470'
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '1.846'
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '149.31'
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("D27").Value = '17.96'
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").Value = '2.089'
$ws.Range("E28").Value = '  -5.47%  '
$ws.Range("D29").Value = '112.79'
$ws.Range("E29").Value = '  -1.94%  '
$ws.Range("D30").Value = '4.672'
$ws.Range("E30").Value = '  -3.96%  '
$ws.Range("D31").Value = '4.642'
$ws.Range("E31").Value = '  -3.03%  '
$ws.Range("D32").Value = '0.09026'
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").Value = '0.05113'
$ws.Range("E33").Value = '  -3.45%  '
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("D35").Value = '1.152'
$ws.Range("E35").Value = '  -6.61%  '
$ws.Range("D36").Value = '0.7301'
$ws.Range("D37").Value = '0.02046'
$ws.Range("E37").Value = '  -1.40%  '
$ws.Range("D38").Value = '2.489'
$ws.Range("E38").Value = '  -5.54%  '
$ws.Range("D39").Value = '3.046'
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("E40").Value = '  -1.97%  '
$ws.Range("D41").Value = '0.5311'
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("D42").Value = '6.581'
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("D43").Value = '115.15'
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("D44").Value = '8.278'
$ws.Range("E44").Value = '  -2.11%  '
$ws.Range("D45").Value = '0.1469'
$ws.Range("E45").Value = '  -2.54%  '
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = '0.4614'
$ws.Range("E47").Value = '  -3.51%  '
$ws.Range("D48").Value = '10.07'
$ws.Range("E48").Value = '  -4.30%  '
$ws.Range("E49").Value = '  -4.21%  '
$ws.Range("D50").Value = '36.63'
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("D51").Value = '64.10'
$ws.Range("E51").Value = '  -4.52%  '
